$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry for 2020-02-04 (row 20)
$ws.Range("A20").Value = 20200204
$ws.Range("B20").Value = "8-9pm"
$ws.Range("C20").Value = "Zihua"
$ws.Range("D20").Value = "Decide features for homework 2."
$ws.Range("F20").Value = "Because we cannot duplicate the features that we have done in homework 1, it took a long time to decide our goals."
$ws.Range("E20").Value = "We decided one feature which is ""Adding more query methods""."
$ws.Range("G20").Value = "Not so good"

# New diary entry for 2020-02-06 (row 21)
$ws.Range("A21").Value = 20200206
$ws.Range("B21").Value = "10am-2pm"
$ws.Range("C21").Value = "Zihua, Wenchia"
$ws.Range("D21").Value = "Decide anther feature and finish writing our report."
$ws.Range("E21").Value = "We decided another one feature which is ""Adding more encryption methods"". The first feature is easy to deal with. But the second one was involved with some C++ native functions which were hard to understand. But we did our best."
$ws.Range("F21").Value = "Because Wenchia had an important interview this morning so we postponed our group activity. We worked hard together and finished the homework in time."
$ws.Range("G21").Value = "Tired"

# Row heights grew to fit the newly-entered wrapped text
$ws.Rows.Item(20).RowHeight = 60.6
$ws.Rows.Item(21).RowHeight = 93.6

# Update the view state to reflect where the user ended up after editing
$ws.Range("E21").Select()
